$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicate fullname: "Phạm Thanh Hà" -> "Phạm Thanh Hà0"
$ws.Range("C2").Value = "Phạm Thanh Hà0"

# Update the active selection to C2 (the edited cell)
$ws.Range("C2").Select()
